$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScoutingData")

$ws.Range("A7").Value = "suday"
# Team Number repeats "1923" from row 3 (B3) - copy it so it stays a text
# shared-string (matching the existing column) instead of being
# auto-detected as a number.
$ws.Range("B3").Copy($ws.Range("B7"))
$ws.Range("C7").Value = "they scuk"
$ws.Range("D7").Value = "its fine"
$ws.Range("E7").Value = "wtv"
